$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# New wellness entries for date 45910 (rows 320-336), matching the two
# format "flavours" already present in the sheet: a blank Localisation
# (col G) cell uses the centered style (like row 318), a filled one
# uses the plain style (like row 319).
$data = @(
    @(320, 45910, "Amir Etien",       76, 9,  10, 4, "Dos",                       3),
    @(321, 45910, "Ilyes Boughanmi",  76, 10, 10, 5, ("Adducteurs" + $nbsp),      0),
    @(322, 45910, "Yanis Berrached",  76, 10, 10, 0, $null,                       0),
    @(323, 45910, "Malik Boussaid",   76, 5,  5,  0, $null,                       10),
    @(324, 45910, "Jeremie Laurent",  76, 8,  8,  0, $null,                       9),
    @(325, 45910, "Kamal Bafounta",   76, 10, 7,  4, "Genou (ménisque)",          6),
    @(326, 45910, "Yoann Martelat",   76, 9,  9,  5, "Genou",                     6),
    @(327, 45910, "Levy Ndoutoume",   76, 8,  7,  4, "Ischio",                    5),
    @(328, 45910, "Emmanuel Valey",   76, 8,  8,  1, ("Adducteur" + $nbsp),       6),
    @(329, 45910, "Ilan Ihaddadene",  76, 9,  9,  0, $null,                       5),
    @(330, 45910, "Karahali Souaré",  76, 8,  8,  7, ("Cheviller" + $nbsp),       3),
    @(331, 45910, "Naim Dhib",        76, 7,  5,  0, $null,                       5),
    @(332, 45910, "Karim Belmahi",    76, 8,  10, 0, $null,                       10),
    @(333, 45910, "Hedi Nasri",       76, 8,  7,  0, $null,                       7),
    @(334, 45910, "Wael Fareh",       76, 7,  7,  1, "Genou",                     6),
    @(335, 45910, "Amir Kherrab",     76, 7,  7,  8, "Cheville",                  6),
    @(336, 45910, "Sofiane Belle",    76, 8,  7,  0, $null,                       7)
)

foreach ($rec in $data) {
    $row = $rec[0]

    # Seed formatting by copying from an existing template row that already
    # has the right style set (blank-G template = row 318, filled-G
    # template = row 319), so we don't mint brand-new style entries.
    if ($null -eq $rec[7]) {
        $ws.Range("A318:I318").Copy()
    } else {
        $ws.Range("A319:I319").Copy()
    }
    $ws.Range("A$row`:I$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $rec[1]
    $ws.Cells.Item($row, 2).Value = $rec[2]
    $ws.Cells.Item($row, 3).Value = $rec[3]
    $ws.Cells.Item($row, 4).Value = $rec[4]
    $ws.Cells.Item($row, 5).Value = $rec[5]
    $ws.Cells.Item($row, 6).Value = $rec[6]
    if ($null -ne $rec[7]) {
        $ws.Cells.Item($row, 7).Value = $rec[7]
    }
    $ws.Cells.Item($row, 8).Value = $rec[8]
    $ws.Cells.Item($row, 9).Formula = "=C$row*D$row"
}

$excel.CutCopyMode = 0

# Mirror the author's final viewport/selection (topLeftCell isn't modelled by
# this headless host, but the active selection round-trips to sheetView).
try { $excel.ActiveWindow.ScrollRow = 315 } catch {}
$ws.Range("L333").Select() | Out-Null
